# (JMT) Added coverage for bl_1s12, 1s16, 1s20, and 1s24
# Update the JMT workspace/user paths used by block bl_1s20 from the old
# "user3/Desktop/jmt_temporary" workspace to the new "user6/jmt_workspace"
# workspace, and update the selected/scrolled cell in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: template path used to build this block
$ws.Range("B1").Value = "template /pub/home/user6/jmt_workspace/blocks/bl_1s20/bl_1s20.tsdl"

# I3: proj_path value
$ws.Range("I3").Value = "/pub/home/user6/jmt_workspace"

# J3: config_path value
$ws.Range("J3").Value = "/pub/home/user6/jmt_workspace/workshop_config.sdl"

# Update the view: scroll so column G is the left-most visible column,
# and select J4 as the active cell.
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1
$ws.Range("J4").Select()
